# Updates the cryptocurrency price (column D) and 1h volume-change
# percentage (column E) figures with freshly scraped values, mirroring
# the automated 'Updated cryptos list ... with GitHub Actions' commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(new Price text (or $null if unchanged), new Volume(1h) text (or $null))
$updates = @{
    2 = @('67.740.62', '  +1.61%  ')
    3 = @('2.622.05', '  +1.67%  ')
    4 = @($null, '  -0.02%  ')
    5 = @('602.54', '  +2.17%  ')
    6 = @('154.63', '  +0.71%  ')
    7 = @($null, '  +0.03%  ')
    8 = @($null, '  +2.18%  ')
    9 = @('2.624.09', '  +1.80%  ')
    10 = @($null, '  +13.49%  ')
    11 = @($null, '  +0.85%  ')
    12 = @('5.24', '  +1.69%  ')
    13 = @('0.355', '  +0.65%  ')
    14 = @('27.92', '  -0.18%  ')
    15 = @('0.0000188', '  +5.53%  ')
    16 = @('3.096.34', $null)
    17 = @('67.647.56', '  +1.74%  ')
    18 = @('2.616.24', '  +1.42%  ')
    19 = @('11.27', '  +0.71%  ')
    20 = @('365.07', '  +3.81%  ')
    21 = @('7.67', '  -0.77%  ')
    22 = @($null, '  -0.07%  ')
    23 = @($null, '  +5.29%  ')
    24 = @($null, '  -0.04%  ')
    25 = @('70.27', '  +4.60%  ')
    26 = @('9.97', '  -3.00%  ')
    27 = @('0.0000105', '  +4.06%  ')
    29 = @('581.70', '  -1.81%  ')
    30 = @($null, '  +0.31%  ')
    31 = @($null, '  +0.19%  ')
    32 = @('7.95', '  -0.08%  ')
    33 = @('1.87', '  +1.18%  ')
    34 = @($null, '  -0.78%  ')
    35 = @('1.00', '  +0.06%  ')
    36 = @($null, '  -1.21%  ')
    37 = @($null, '  +0.40%  ')
    38 = @('158.34', '  +3.41%  ')
    39 = @($null, '  +2.25%  ')
    40 = @($null, '  +1.26%  ')
    41 = @('5.40', '  -0.11%  ')
    42 = @('1.85', '  +4.36%  ')
    43 = @('2.68', '  +3.66%  ')
    44 = @('41.15', '  -0.67%  ')
    45 = @('16.43', '  +0.12%  ')
    46 = @('0.999', '  -0.02%  ')
    47 = @('157.25', '  +1.30%  ')
    48 = @('0.0₆0292', '  -4.33%  ')
    49 = @($null, '  +1.05%  ')
    50 = @('21.04', '  +0.55%  ')
    51 = @($null, '  +2.17%  ')
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $newPrice = $pair[0]
    $newVolume = $pair[1]

    if ($newPrice -ne $null) {
        $cell = $ws.Cells.Item($row, 4)  # column D = Price
        # Column D mixes numeric-looking text ('602.54') with
        # genuinely non-numeric text ('67.740.62', '0.0₆0292').
        # Force text entry for the numeric-looking ones (leading
        # apostrophe) so Excel keeps them as strings, like the
        # original inline-string cells, then restore the default
        # 'Normal' style so no stray number-format/style sticks.
        if ($newPrice -match '^-?\d+(\.\d+)?$') {
            $cell.Value = "'" + $newPrice
            $cell.Style = 'Normal'
        } else {
            $cell.Value = $newPrice
        }
    }

    if ($newVolume -ne $null) {
        $ws.Cells.Item($row, 5).Value = $newVolume  # column E = Volume(1h)
    }
}
